$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.617.35"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.621.73"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = $ws.Range("D49").Style
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'604.30"
$ws.Range("D5").Style = $ws.Range("D49").Style
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'197.02"
$ws.Range("D6").Style = $ws.Range("D49").Style
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = $ws.Range("D49").Style
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "'0.648"
$ws.Range("D10").Style = $ws.Range("D49").Style
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'53.65"
$ws.Range("D11").Style = $ws.Range("D49").Style
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "4.195.44"
$ws.Range("E14").Value = "  +2.25%  "
$ws.Range("D15").Value = "'600.29"
$ws.Range("D15").Style = $ws.Range("D49").Style
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "'13.01"
$ws.Range("D16").Style = $ws.Range("D49").Style
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "70.644.10"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "3.620.70"
$ws.Range("E18").Value = "  +2.19%  "
$ws.Range("D19").Value = "'19.07"
$ws.Range("D19").Style = $ws.Range("D49").Style
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "'101.78"
$ws.Range("D24").Style = $ws.Range("D49").Style
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("E25").Value = "  +0.19%  "
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.47%  "
$ws.Range("D30").Value = "'4.69"
$ws.Range("D30").Style = $ws.Range("D49").Style
$ws.Range("E30").Value = "  +7.69%  "
$ws.Range("E31").Value = "  +1.28%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").Value = "'63.44"
$ws.Range("D34").Style = $ws.Range("D49").Style
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "0.0₃0888"
$ws.Range("E35").Value = "  +5.87%  "
$ws.Range("D36").Value = "3.910.19"
$ws.Range("E36").Value = "  +3.26%  "
$ws.Range("D37").Value = "'540.59"
$ws.Range("D37").Style = $ws.Range("D49").Style
$ws.Range("E37").Value = "  +9.53%  "
$ws.Range("D38").Value = "'3.12"
$ws.Range("D38").Style = $ws.Range("D49").Style
$ws.Range("E38").Value = "  +1.45%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").Value = "'37.04"
$ws.Range("D40").Style = $ws.Range("D49").Style
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("E41").Value = "  -1.33%  "
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "'0.0460"
$ws.Range("D44").Style = $ws.Range("D49").Style
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'3.41"
$ws.Range("D45").Style = $ws.Range("D49").Style
$ws.Range("E45").Value = "  +3.72%  "
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'0.141"
$ws.Range("D47").Style = $ws.Range("D49").Style
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "'8.61"
$ws.Range("D48").Style = $ws.Range("D49").Style
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("E51").Value = "  +2.38%  "
